# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
#
# Note: several "Price" values are plain decimal-looking strings
# (e.g. "25.42", "1.00") that must stay TEXT (matching the workbook's
# original inlineStr cells), not be coerced into numbers by Excel's
# normal text auto-detection. Writing them through .Value directly
# would turn "25.42" into the number 25.42. To avoid that we instead
# write them as a string-literal formula (="25.42") and then do a
# Copy / PasteSpecial-values-only pass, which freezes the formula's
# cached text result back down to a literal cell without Excel
# re-parsing the string as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace("""", """""")
    $ws.Range($range).Formula = "=""" + $escaped + """"
    $ws.Range($range).Copy()
    $ws.Range($range).PasteSpecial(-4163)
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "67.203.13"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "2.490.02"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 (BNB)
Set-TextValue "D5" "584.50"
$ws.Range("E5").Value = "  -0.18%  "

# Row 6 (Solana)
Set-TextValue "D6" "172.72"
$ws.Range("E6").Value = "  +2.78%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 (XRP)
$ws.Range("E8").Value = "  -0.86%  "

# Row 9 (LidoStakedEther)
Set-TextValue "D9" "2.488.99"
$ws.Range("E9").Value = "  +0.10%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +0.59%  "

# Row 11 (TRON)
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 (Toncoin)
$ws.Range("E12").Value = "  -0.38%  "

# Row 13 (Cardano)
$ws.Range("E13").Value = "  -1.78%  "

# Row 15 (Avalanche)
Set-TextValue "D15" "25.42"
$ws.Range("E15").Value = "  -2.02%  "

# Row 16 (WrappedBTC)
Set-TextValue "D16" "67.207.73"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17 (ShibaInu)
$ws.Range("E17").Value = "  -1.76%  "

# Row 18 (WrappedEther)
Set-TextValue "D18" "2.496.00"
$ws.Range("E18").Value = "  -0.92%  "

# Row 19 (Chainlink)
$ws.Range("E19").Value = "  -6.13%  "

# Row 20 (Uniswap)
Set-TextValue "D20" "7.42"
$ws.Range("E20").Value = "  -4.81%  "

# Row 21 (BitcoinCash)
Set-TextValue "D21" "349.36"
$ws.Range("E21").Value = "  -3.05%  "

# Row 22 (Polkadot)
$ws.Range("E22").Value = "  -0.56%  "

# Row 23 (Dai)
$ws.Range("E23").Value = "  -0.15%  "

# Row 24 (Litecoin)
Set-TextValue "D24" "68.56"
$ws.Range("E24").Value = "  -3.31%  "

# Row 25 (NEARProtocol)
$ws.Range("E25").Value = "  -4.24%  "

# Row 26 (SuiNetwork)
$ws.Range("E26").Value = "  -2.98%  "

# Row 27 (Aptos)
Set-TextValue "D27" "9.27"
$ws.Range("E27").Value = "  -2.02%  "

# Row 28 (Binance-PegBSC-USD)
Set-TextValue "D28" "1.00"

# Row 29 (WrappedeETH)
$ws.Range("E29").Value = "  +0.24%  "

# Row 30 (PEPE)
$ws.Range("E30").Value = "  -3.47%  "

# Row 31 (Bittensor)
Set-TextValue "D31" "509.16"
$ws.Range("E31").Value = "  +0.34%  "

# Row 32 (InternetComputer(DFINITY))
$ws.Range("E32").Value = "  -3.58%  "

# Row 33 (Fetch.AI)
Set-TextValue "D33" "1.24"
$ws.Range("E33").Value = "  -3.21%  "

# Row 34 (PancakeSwap)
$ws.Range("E34").Value = "  -3.81%  "

# Row 35 (FirstDigitalUSD)
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36 (Monero)
Set-TextValue "D36" "159.81"
$ws.Range("E36").Value = "  +0.75%  "

# Row 37 (Kaspa)
$ws.Range("E37").Value = "  -7.46%  "

# Row 38 (WhiteBITCoin)
Set-TextValue "D38" "18.70"
$ws.Range("E38").Value = "  +0.77%  "

# Row 39 (EthereumClassic)
Set-TextValue "D39" "18.23"
$ws.Range("E39").Value = "  -4.33%  "

# Row 40 (ImmutableX)
$ws.Range("E40").Value = "  -5.71%  "

# Row 41 (Stacks)
$ws.Range("E41").Value = "  -2.58%  "

# Row 42 (USDe)
$ws.Range("E42").Value = "  -0.06%  "

# Row 43 (RenderToken)
$ws.Range("E43").Value = "  -2.66%  "

# Row 44 (PolygonEcosystemToken)
$ws.Range("E44").Value = "  -2.06%  "

# Row 45 (dogwifhat)
Set-TextValue "D45" "2.36"
$ws.Range("E45").Value = "  -4.31%  "

# Row 46 (OKB)
Set-TextValue "D46" "38.59"
$ws.Range("E46").Value = "  -2.09%  "

# Row 47 (Aave)
Set-TextValue "D47" "142.74"
$ws.Range("E47").Value = "  +0.49%  "

# Row 48 (ARBITRUM)
$ws.Range("E48").Value = "  -4.66%  "

# Row 49 (Filecoin)
$ws.Range("E49").Value = "  -4.25%  "

# Row 50 (BabyDogeCoin)
Set-TextValue "D50" "0.0₆0250"
$ws.Range("E50").Value = "  -6.48%  "

# Row 51 (Cronos)
$ws.Range("E51").Value = "  -0.89%  "
